# Apply "User data with filename fix" changes:
#  - Notes sheet: fix "Units of measure" value
#  - Data sheet: populate the previously-empty data rows

$wb = $excel.ActiveWorkbook

# --- Notes sheet -----------------------------------------------------
$notes = $wb.Worksheets.Item("Notes")
$notes.Range("A3").Value = "Units of measure: constant 2015 US$"

# --- Data sheet --------------------------------------------------------
$data = $wb.Worksheets.Item("Data")

$rows = @(
    @("bilateral-unspecified", "Bilateral, unspecified", 2014, 0),
    @("CV", "Cape Verde", 2014, 330156.66),
    @("CV", "Cape Verde", 2015, 500000),
    @("GW", "Guinea-Bissau", 2014, 462219.32),
    @("ID", "Indonesia", 2014, 0),
    @("MY", "Malaysia", 2015, 1000000),
    @("MM", "Myanmar", 2015, 440000),
    @("NP", "Nepal", 2015, 500000),
    @("KP", "North Korea", 2014, 66031.33),
    @("north-of-sahara", "North of Sahara, regional", 2015, 60000),
    @("ST", "Sao Tome & Principe", 2015, 0),
    @("south-of-sahara", "South of Sahara, regional", 2014, 1320626.62),
    @("TH", "Thailand", 2015, 500000),
    @("TO", "Tonga", 2014, 0),
    @("VU", "Vanuatu", 2015, 1000000)
)

$r = 2
foreach ($row in $rows) {
    $data.Cells.Item($r, 1).Value = $row[0]
    $data.Cells.Item($r, 2).Value = $row[1]
    $data.Cells.Item($r, 3).Value = $row[2]
    $data.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
